$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.037522333333333
$ws.Range("H2").Value = 9.112567
$ws.Range("I2").Value = 0.1153015356242242
$ws.Range("J2").Value = 0.1153015356242242
$ws.Range("M2").Value = 0.7521946666666667
$ws.Range("N2").Value = 2.256584
$ws.Range("O2").Value = 0.07361670343069449
$ws.Range("P2").Value = 0.0736167034306945
$ws.Range("Q2").Value = 2.284808099014223
$ws.Range("R2").Value = 20.563272891128
$ws.Range("S2").Value = 0.008488118953152168
$ws.Range("T2").Value = 0.00848811895315217

$ws.Range("G3").Value = 3.037522333333333
$ws.Range("H3").Value = 9.112567
$ws.Range("I3").Value = 0.1153015356242242
$ws.Range("J3").Value = 0.1153015356242242
$ws.Range("M3").Value = 7.059280333333334
$ws.Range("N3").Value = 21.177841
$ws.Range("O3").Value = 0.6908862423022597
$ws.Range("P3").Value = 0.6908862423022598
$ws.Range("Q3").Value = 21.44272166976078
$ws.Range("R3").Value = 192.984495027847
$ws.Range("S3").Value = 0.07966024467910039
$ws.Range("T3").Value = 0.0796602446791004

$ws.Range("G4").Value = 3.037522333333333
$ws.Range("H4").Value = 9.112567
$ws.Range("I4").Value = 0.1153015356242242
$ws.Range("J4").Value = 0.1153015356242242
$ws.Range("M4").Value = 2.406242333333334
$ws.Range("N4").Value = 7.218727
$ws.Range("O4").Value = 0.2354970542670457
$ws.Range("P4").Value = 0.2354970542670457
$ws.Range("Q4").Value = 7.309014826912112
$ws.Range("R4").Value = 65.78113344220901
$ws.Range("S4").Value = 0.02715317199197163
$ws.Range("T4").Value = 0.02715317199197163

$ws.Range("G5").Value = 8.588082333333332
$ws.Range("H5").Value = 25.764247
$ws.Range("I5").Value = 0.325995654495798
$ws.Range("J5").Value = 0.325995654495798
$ws.Range("M5").Value = 0.7521946666666667
$ws.Range("N5").Value = 2.256584
$ws.Range("O5").Value = 0.07361670343069449
$ws.Range("P5").Value = 0.0736167034306945
$ws.Range("Q5").Value = 6.459909728027555
$ws.Range("R5").Value = 58.139187552248
$ws.Range("S5").Value = 0.02399872541671231
$ws.Range("T5").Value = 0.02399872541671231

$ws.Range("G6").Value = 8.588082333333332
$ws.Range("H6").Value = 25.764247
$ws.Range("I6").Value = 0.325995654495798
$ws.Range("J6").Value = 0.325995654495798
$ws.Range("M6").Value = 7.059280333333334
$ws.Range("N6").Value = 21.177841
$ws.Range("O6").Value = 0.6908862423022597
$ws.Range("P6").Value = 0.6908862423022598
$ws.Range("Q6").Value = 60.62568071674744
$ws.Range("R6").Value = 545.6311264507269
$ws.Range("S6").Value = 0.2252259127414677
$ws.Range("T6").Value = 0.2252259127414677

$ws.Range("G7").Value = 8.588082333333332
$ws.Range("H7").Value = 25.764247
$ws.Range("I7").Value = 0.325995654495798
$ws.Range("J7").Value = 0.325995654495798
$ws.Range("M7").Value = 2.406242333333334
$ws.Range("N7").Value = 7.218727
$ws.Range("O7").Value = 0.2354970542670457
$ws.Range("P7").Value = 0.2354970542670457
$ws.Range("Q7").Value = 20.66500727261878
$ws.Range("R7").Value = 185.985065453569
$ws.Range("S7").Value = 0.07677101633761803
$ws.Range("T7").Value = 0.07677101633761803

$ws.Range("G8").Value = 14.71855733333333
$ws.Range("H8").Value = 44.155672
$ws.Range("I8").Value = 0.5587028098799778
$ws.Range("J8").Value = 0.5587028098799777
$ws.Range("M8").Value = 0.7521946666666667
$ws.Range("N8").Value = 2.256584
$ws.Range("O8").Value = 0.07361670343069449
$ws.Range("P8").Value = 0.0736167034306945
$ws.Range("Q8").Value = 11.07122032716089
$ws.Range("R8").Value = 99.64098294444801
$ws.Range("S8").Value = 0.04112985906083001
$ws.Range("T8").Value = 0.04112985906083001

$ws.Range("G9").Value = 14.71855733333333
$ws.Range("H9").Value = 44.155672
$ws.Range("I9").Value = 0.5587028098799778
$ws.Range("J9").Value = 0.5587028098799777
$ws.Range("M9").Value = 7.059280333333334
$ws.Range("N9").Value = 21.177841
$ws.Range("O9").Value = 0.6908862423022597
$ws.Range("P9").Value = 0.6908862423022598
$ws.Range("Q9").Value = 103.9024223182391
$ws.Range("R9").Value = 935.1218008641521
$ws.Range("S9").Value = 0.3860000848816917
$ws.Range("T9").Value = 0.3860000848816916

$ws.Range("G10").Value = 14.71855733333333
$ws.Range("H10").Value = 44.155672
$ws.Range("I10").Value = 0.5587028098799778
$ws.Range("J10").Value = 0.5587028098799777
$ws.Range("M10").Value = 2.406242333333334
$ws.Range("N10").Value = 7.218727
$ws.Range("O10").Value = 0.2354970542670457
$ws.Range("P10").Value = 0.2354970542670457
$ws.Range("Q10").Value = 35.41641574106045
$ws.Range("R10").Value = 318.747741669544
$ws.Range("S10").Value = 0.131572865937456
$ws.Range("T10").Value = 0.131572865937456
